$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; temporarily unprotect so the cells below can be edited
$ws.Unprotect()

# Update the confidential disclaimer text with the new "as of" date
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-21 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-10
$ws.Range("D2").Value = 0.09211912044072743
$ws.Range("E2").Value = -0.005078075409419736

$ws.Range("D3").Value = 0.1064495703270682
$ws.Range("E3").Value = -0.00392251523746312

$ws.Range("D4").Value = 0.119818584430478
$ws.Range("E4").Value = -0.00148879485973985

$ws.Range("D5").Value = 0.140779762238298
$ws.Range("E5").Value = 0.001113851539501942

$ws.Range("D6").Value = 0.1374818579645535
$ws.Range("E6").Value = -0.001515569027280228

$ws.Range("D7").Value = 0.1463293072783039
$ws.Range("E7").Value = 0.00425491679273815

$ws.Range("D8").Value = 0.127328885262845
$ws.Range("E8").Value = -0.002380243975007357

$ws.Range("D9").Value = 0.1296929120577259
$ws.Range("E9").Value = 0.004888844911946411

$ws.Range("E10").Value = -0.0001616849372765294

# Restore sheet protection (same flags as the original workbook)
$ws.Protect("D382", $false, $true, $true, $true, $false, $false, $false, $false, $false, $false, $false, $false, $false, $false, $false)
